# Auto-generated from diff: updates NATMI Pdgfb-Pdgfrb LR-pair TPM recalculation values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 46.11811066666667
$ws.Range("H2").Value = 138.354332
$ws.Range("I2").Value = 0.95896098489411
$ws.Range("J2").Value = 0.9589609848941099
$ws.Range("M2").Value = 8.488196666666667
$ws.Range("N2").Value = 25.46459
$ws.Range("O2").Value = 0.04138402976425696
$ws.Range("P2").Value = 0.04138402976425696
$ws.Range("Q2").Value = 391.4595932337645
$ws.Range("R2").Value = 3523.13633910388
$ws.Range("S2").Value = 0.03968566994161901
$ws.Range("T2").Value = 0.03968566994161901

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 46.11811066666667
$ws.Range("H3").Value = 138.354332
$ws.Range("I3").Value = 0.95896098489411
$ws.Range("J3").Value = 0.9589609848941099
$ws.Range("O3").Value = 0.3297460182766552
$ws.Range("P3").Value = 0.3297460182766552
$ws.Range("Q3").Value = 3119.131774270087
$ws.Range("R3").Value = 28072.18596843078
$ws.Range("S3").Value = 0.3162135664514925
$ws.Range("T3").Value = 0.3162135664514924

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 46.11811066666667
$ws.Range("H4").Value = 138.354332
$ws.Range("I4").Value = 0.95896098489411
$ws.Range("J4").Value = 0.9589609848941099
$ws.Range("O4").Value = 0.6288699519590879
$ws.Range("P4").Value = 0.6288699519590879
$ws.Range("Q4").Value = 5948.603289558397
$ws.Range("R4").Value = 53537.42960602557
$ws.Range("S4").Value = 0.6030617485009986
$ws.Range("T4").Value = 0.6030617485009985

# Row 5
$ws.Range("I5").Value = 0.002799731840346333
$ws.Range("J5").Value = 0.002799731840346333
$ws.Range("M5").Value = 8.488196666666667
$ws.Range("N5").Value = 25.46459
$ws.Range("O5").Value = 0.04138402976425696
$ws.Range("P5").Value = 0.04138402976425696
$ws.Range("Q5").Value = 1.142884751986667
$ws.Range("R5").Value = 10.28596276788
$ws.Range("S5").Value = 0.0001158641858128306
$ws.Range("T5").Value = 0.0001158641858128306

# Row 6
$ws.Range("I6").Value = 0.002799731840346333
$ws.Range("J6").Value = 0.002799731840346333
$ws.Range("O6").Value = 0.3297460182766552
$ws.Range("P6").Value = 0.3297460182766552
$ws.Range("Q6").Value = 9.106452379419999
$ws.Range("R6").Value = 81.95807141477999
$ws.Range("S6").Value = 0.0009232004265965753
$ws.Range("T6").Value = 0.0009232004265965753

# Row 7
$ws.Range("I7").Value = 0.002799731840346333
$ws.Range("J7").Value = 0.002799731840346333
$ws.Range("O7").Value = 0.6288699519590879
$ws.Range("P7").Value = 0.6288699519590879
$ws.Range("S7").Value = 0.001760667227936927
$ws.Range("T7").Value = 0.001760667227936927

# Row 8
$ws.Range("I8").Value = 0.0382392832655437
$ws.Range("J8").Value = 0.0382392832655437
$ws.Range("M8").Value = 8.488196666666667
$ws.Range("N8").Value = 25.46459
$ws.Range("O8").Value = 0.04138402976425696
$ws.Range("P8").Value = 0.04138402976425696
$ws.Range("Q8").Value = 15.60974274082
$ws.Range("R8").Value = 140.48768466738
$ws.Range("S8").Value = 0.001582495636825113
$ws.Range("T8").Value = 0.001582495636825113

# Row 9
$ws.Range("I9").Value = 0.0382392832655437
$ws.Range("J9").Value = 0.0382392832655437
$ws.Range("O9").Value = 0.3297460182766552
$ws.Range("P9").Value = 0.3297460182766552
$ws.Range("S9").Value = 0.01260925139856617
$ws.Range("T9").Value = 0.01260925139856617

# Row 10
$ws.Range("I10").Value = 0.0382392832655437
$ws.Range("J10").Value = 0.0382392832655437
$ws.Range("O10").Value = 0.6288699519590879
$ws.Range("P10").Value = 0.6288699519590879
$ws.Range("S10").Value = 0.02404753623015242
$ws.Range("T10").Value = 0.02404753623015242
